# cafe_management/src/data/dataTemplate.xlsx
# "Users" sheet: replace the old user1/pass1 credentials with u1/p1, and
# append a new 4th user row (u04) with numeric placeholder 1/1 values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Row 2 previously held user1 / pass1 -> now u1 / p1
$ws.Range("B2").Value = "u1"
$ws.Range("C2").Value = "p1"

# New row 5: Id = u04, with numeric 1 / 1 in the Username/Password columns
$ws.Range("A5").Value = "u04"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 1

# Carry the existing row-style (vertical-center + wrap text) down onto the
# new row's B5:C5 cells, matching the formatting used by rows 2-4.
$ws.Range("B4:C4").Copy()
$ws.Range("B5:C5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Update the active selection shown in the sheet view (D6 -> C6)
$ws.Activate() | Out-Null
$ws.Range("C6").Select() | Out-Null
